# Daily attendance processing - swap the order of names in the
# "Recorded By" column (G) from "System, dnasr281@gmail.com" to
# "dnasr281@gmail.com, System" for every row where it occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$colG = $ws.Range("G1:G$lastRow")

$colG.Replace($oldValue, $newValue)
